$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "theta_threshold_range" row (row 5: Min=11.2, Max=176) is removed entirely.
# Deleting the row shifts the following row ("pie_threshold_range") up into row 5,
# and the engine also drops the now-unused shared string automatically.
$ws.Rows("5").Delete()

# The row that is now row 5 ("pie_threshold_range") gets new Min/Max values.
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# Update the saved selection/active cell.
$ws.Range("C4").Select()

# Configure the print page setup (paper size + portrait orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
